# Apply the same text edit that Chandana Ch made to the "ChitFund DApp"
# slide (sldId 281, the 11th slide): the bullet that talks about creating
# separate underwriter functions keeps its text, but the following bullet
# ("Have couple of compilation issues") is removed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item("Subtitle 2")
$tr = $shape.TextFrame.TextRange

# Paragraph 7 = "Create separate functions to verify the underwriters,
# assign an underwriter and decision making". Re-stamp the trailing
# "decision making" portion so it becomes its own run (matches the
# split produced by the original edit).
$para7 = $tr.Paragraphs(7)
$marker = "decision making"
$relIdx = $para7.Text.IndexOf($marker)
$tailStart = $para7.Start + $relIdx
$tail = $tr.Characters($tailStart, $marker.Length)
$tail.Text = $marker

# Paragraph 8 = "Have couple of compilation issues" bullet - remove it
# entirely.
$para8 = $tr.Paragraphs(8)
$para8.Delete()
